$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the user name values in column A (rows 4, 6, 7) to reflect the
# renumbered manufacturer/authorised-rep test fixtures (126 -> 196).
$ws.Range("A4").Value = "AuthorisedRep1961_AT"
$ws.Range("A6").Value = "Manufacturer196_NU"
$ws.Range("A7").Value = "AuthorisedRep196_NU"

# Keep the "Manufacturer126_AT" fixture name in sync with the rest (A3
# already references this shared string, so updating its text updates it
# everywhere it is used).
$ws.Range("A3").Value = "Manufacturer196_AT"

# Move the active selection from C8 to B5, as recorded in the sheet view.
$ws.Activate()
$ws.Range("B5").Select()
